$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new suspension rod parameters in rows 18-21
$ws.Range("A18").Value = "Pull_rod_Offset (mm)"
$ws.Range("B18").Value = 30

$ws.Range("A19").Value = "Push_rod_Offset (mm)"
$ws.Range("B19").Value = 30

$ws.Range("A20").Value = "Suspension_Rod_Rint (mm)"
$ws.Range("B20").Value = 9

$ws.Range("A21").Value = "Suspension_Rod_Rext (mm)"
$ws.Range("B21").Value = 10

# Match number formats from existing rows (A column = text "@", B column = number "0.00")
$ws.Range("A18:A21").NumberFormat = $ws.Range("A17").NumberFormat
$ws.Range("B18:B21").NumberFormat = $ws.Range("B17").NumberFormat

# Update selection to match the diff (activeCell F22)
$ws.Range("F22").Select()
